$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.778.52"
$ws.Range("E2").Value = "  -3.72%  "
$ws.Range("D3").Value = "3.498.58"
$ws.Range("E3").Value = "  -3.03%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.59%  "
$ws.Range("D7").Value = "3.495.88"
$ws.Range("E7").Value = "  -3.22%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.384"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.69%  "
$ws.Range("D13").Value = "4.100.80"
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000179"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.66%  "
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").Value = "3.506.55"
$ws.Range("E17").Value = "  -2.76%  "
$ws.Range("D18").Value = "63.968.26"
$ws.Range("E18").Value = "  -3.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.578"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("D24").Value = "3.642.36"
$ws.Range("E24").Value = "  -2.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.07%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000112"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.28%  "
$ws.Range("D33").Value = "3.504.78"
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.47%  "
$ws.Range("E36").Value = "  -4.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.15%  "
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0805"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.812"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.95%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.39%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.24%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.23%  "
$ws.Range("D50").Value = "2.441.61"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.901"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.25%  "
